$d = $word.ActiveDocument

# 1. "Power Query (Ferramenta ETL)" bullet (Modulo 2) -> prepend "Automatizando tarefas com "
$rng1 = $d.Content
$rng1.Find.Execute("Power Query (Ferramenta ETL)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Collapse(1)
$pos1 = $rng1.Start
$ins1 = $d.Range($pos1, $pos1)
$ins1.Text = "Automatizando tarefas com "

# 2. "Modelagem e transformação de dados no Power Query" -> append " (ETL)"
$rng2 = $d.Content
$rng2.Find.Execute("Modelagem e transformação de dados no Power Query", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$pos2 = $rng2.Start
$ins2 = $d.Range($pos2, $pos2)
$ins2.Text = " (ETL)"
